$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feeds")

# Clear redundant fallback URL cells in column D that duplicated column B/C values
$ws.Range("D8").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("D23").ClearContents()
$ws.Range("D24").ClearContents()
$ws.Range("D28").ClearContents()
$ws.Range("D29").ClearContents()

# Add the two new Slovenian feed URLs on row 25 (SI)
$ws.Range("B25").Value = "https://zurnal24.si/feeds/latest"
$ws.Range("C25").Value = "https://24ur.com/rss"

# Update the active selection to C25 to match the saved view state
$ws.Range("C25").Select()
